# Actualizacion tablas al guardar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename "tartaleta" -> "tartaleta durazno" in the product table (A5)
$ws.Range("A5").Value = "tartaleta durazno"

# Reflect the last active selection being on A5 when the file was saved
$ws.Range("A5").Select()
